# Rename the customer labels in column A (rows 2-5) of Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 1).Value = "Alpha"
$ws.Cells.Item(3, 1).Value = "Gamma"
$ws.Cells.Item(4, 1).Value = "Beta"
$ws.Cells.Item(5, 1).Value = "Sigma"

# Update the active selection to match the saved view state.
$ws.Range("F10").Select()
